# Updated cryptos list - apply new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.147.21"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.930.13"
$ws.Range("E3").Value = "  +4.24%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "353.03"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "113.09"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("D10").Value = "39.59"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "0.0881"
$ws.Range("E11").Value = "  +4.24%  "
$ws.Range("D13").Value = "20.03"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "3.396.80"
$ws.Range("E15").Value = "  +4.41%  "
$ws.Range("D16").Value = "2.922.17"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "0.983"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "52.218.57"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").Value = "'7.60"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "3.29"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "14.19"
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "71.12"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "269.07"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  +9.75%  "
$ws.Range("D27").Value = "27.01"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'7.00"
$ws.Range("E29").Value = "  +13.65%  "
$ws.Range("D30").Value = "10.65"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "0.104"
$ws.Range("E31").Value = "  +15.06%  "
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "37.17"
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("D34").Value = "6.04"
$ws.Range("E34").Value = "  +5.91%  "
$ws.Range("D35").Value = "'53.00"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("D39").Value = "18.67"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "23.16"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("D45").Value = "2.193.40"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").Value = "3.53"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "111.64"
$ws.Range("E48").Value = "  -7.59%  "
$ws.Range("E49").Value = "  +10.00%  "
$ws.Range("D50").Value = "0.0346"
$ws.Range("E50").Value = "  +7.86%  "
$ws.Range("D51").Value = "0.956"
$ws.Range("E51").Value = "  -5.90%  "
